$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the example clinical dates (row 2)
$ws.Range("C2").Value = 42537
$ws.Range("D2").Value = 42761
$ws.Range("F2").Value = 42956
$ws.Range("G2").Value = 43116

# Reposition the window and the active selection in the frozen-pane view
$win = $excel.ActiveWindow
$win.Left = -31720
$win.Top = 6180

$ws.Range("G5").Select()
